$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 57695080
$ws.Range("I76").Value = 68184560
$ws.Range("J76").Value = 2950
$ws.Range("K76").Value = 68184560
$ws.Range("L76").Value = 2950
$ws.Range("M76").Value = -68184245
$ws.Range("N76").Value = -3580

$ws.Range("H79").Value = 57695080
$ws.Range("I79").Value = 68184560
$ws.Range("J79").Value = 2950
$ws.Range("K79").Value = 68184560
$ws.Range("L79").Value = 2950
$ws.Range("M79").Value = -68183468
$ws.Range("N79").Value = -5134

$ws.Range("H98").Value = 1893.2941
$ws.Range("J98").Value = 904
$ws.Range("L98").Value = 904
$ws.Range("N98").Value = -3900

$ws.Range("H113").Value = 3316.4614
$ws.Range("I113").Value = 3510.5715
$ws.Range("J113").Value = 3090
$ws.Range("K113").Value = 3510.5715
$ws.Range("L113").Value = 3090
$ws.Range("M113").Value = -256.5715
$ws.Range("N113").Value = -9598

$ws.Range("H122").Value = 1893.2941
$ws.Range("J122").Value = 904
$ws.Range("L122").Value = 2712
$ws.Range("N122").Value = -7612

$ws.Range("H137").Value = 995.9524
$ws.Range("I137").Value = 806.64703
$ws.Range("K137").Value = 2419.94109
$ws.Range("M137").Value = 130.0589100000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1253.2909
$ws.Range("I61").Value = 1301.5111
$ws.Range("K61").Value = 1301.5111
$ws.Range("M61").Value = -1089.5111

$ws.Range("H63").Value = 1992.5625
$ws.Range("I63").Value = 1927.5834
$ws.Range("K63").Value = 1927.5834
$ws.Range("M63").Value = -1241.5834

$ws.Range("H66").Value = 1992.5625
$ws.Range("I66").Value = 1927.5834
$ws.Range("K66").Value = 9637.916999999999
$ws.Range("M66").Value = -6205.916999999999

$ws.Range("H136").Value = 1253.2909
$ws.Range("I136").Value = 1301.5111
$ws.Range("K136").Value = 3904.5333
$ws.Range("M136").Value = -1354.5333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 52633730
$ws.Range("I99").Value = 76925320
$ws.Range("J99").Value = 1933.3334
$ws.Range("K99").Value = 76925320
$ws.Range("L99").Value = 1933.3334
$ws.Range("M99").Value = -76923822
$ws.Range("N99").Value = -4929.3334

$ws.Range("H105").Value = 3635.093
$ws.Range("I105").Value = 3567.92
$ws.Range("J105").Value = 3728.389
$ws.Range("K105").Value = 3567.92
$ws.Range("L105").Value = 3728.389
$ws.Range("M105").Value = -1820.92
$ws.Range("N105").Value = -7222.389

$ws.Range("H134").Value = 37449.17
$ws.Range("I134").Value = 2691.0476
$ws.Range("J134").Value = 128689.25
$ws.Range("K134").Value = 8073.1428
$ws.Range("L134").Value = 386067.75
$ws.Range("M134").Value = -5538.1428
$ws.Range("N134").Value = -391137.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 11918.556
$ws.Range("I74").Value = 5750
$ws.Range("J74").Value = 13681
$ws.Range("K74").Value = 5750
$ws.Range("L74").Value = 13681
$ws.Range("M74").Value = -4876
$ws.Range("N74").Value = -15429

$ws.Range("H77").Value = 11918.556
$ws.Range("I77").Value = 5750
$ws.Range("J77").Value = 13681
$ws.Range("K77").Value = 17250
$ws.Range("L77").Value = 41043
$ws.Range("M77").Value = -12882
$ws.Range("N77").Value = -49779

$ws.Range("H94").Value = 9760.857
$ws.Range("J94").Value = 10052.333
$ws.Range("L94").Value = 10052.333
$ws.Range("N94").Value = -10954.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2941306.2
$ws.Range("I2").Value = 4902149
$ws.Range("J2").Value = 42.25
$ws.Range("K2").Value = 29412894
$ws.Range("L2").Value = 253.5
$ws.Range("M2").Value = -29412781
$ws.Range("N2").Value = -479.5

$ws.Range("H92").Value = 38217.125
$ws.Range("I92").Value = 100467.336
$ws.Range("J92").Value = 867
$ws.Range("K92").Value = 301402.008
$ws.Range("L92").Value = 2601
$ws.Range("M92").Value = -300154.008
$ws.Range("N92").Value = -5097

$ws.Range("H131").Value = 20876818
$ws.Range("I131").Value = 125252580
$ws.Range("J131").Value = 1667.45
$ws.Range("K131").Value = 375757740
$ws.Range("L131").Value = 5002.35
$ws.Range("M131").Value = -375752700
$ws.Range("N131").Value = -15082.35

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4735.5
$ws.Range("I70").Value = 4629.7
$ws.Range("K70").Value = 4629.7
$ws.Range("M70").Value = -4359.7

$ws.Range("H73").Value = 4735.5
$ws.Range("I73").Value = 4629.7
$ws.Range("K73").Value = 4629.7
$ws.Range("M73").Value = -3693.7

$ws.Range("H80").Value = 3861.5386
$ws.Range("I80").Value = 4070
$ws.Range("J80").Value = 3166.6667
$ws.Range("K80").Value = 4070
$ws.Range("L80").Value = 3166.6667
$ws.Range("M80").Value = -3072
$ws.Range("N80").Value = -5162.6667

$ws.Range("H83").Value = 3861.5386
$ws.Range("I83").Value = 4070
$ws.Range("J83").Value = 3166.6667
$ws.Range("K83").Value = 20350
$ws.Range("L83").Value = 15833.3335
$ws.Range("M83").Value = -15358
$ws.Range("N83").Value = -25817.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1498.6666
$ws.Range("I7").Value = 1500.7368
$ws.Range("J7").Value = 1493.75
$ws.Range("K7").Value = 1500.7368
$ws.Range("L7").Value = 1493.75
$ws.Range("M7").Value = -1388.7368
$ws.Range("N7").Value = -1717.75

$ws.Range("H122").Value = 2430.15
$ws.Range("I122").Value = 2437
$ws.Range("J122").Value = 2414.1667
$ws.Range("K122").Value = 7311
$ws.Range("L122").Value = 7242.500100000001
$ws.Range("M122").Value = -4861
$ws.Range("N122").Value = -12142.5001

$ws.Range("H126").Value = 1498.6666
$ws.Range("I126").Value = 1500.7368
$ws.Range("J126").Value = 1493.75
$ws.Range("K126").Value = 4502.2104
$ws.Range("L126").Value = 4481.25
$ws.Range("M126").Value = -2032.2104
$ws.Range("N126").Value = -9421.25

$ws.Range("H136").Value = 4095.625
$ws.Range("I136").Value = 2530
$ws.Range("J136").Value = 7540
$ws.Range("K136").Value = 7590
$ws.Range("L136").Value = 22620
$ws.Range("M136").Value = -5040
$ws.Range("N136").Value = -27720
